$d = $word.ActiveDocument

# Helper: split the run containing absolute document position $pos into two
# runs (both keeping identical formatting) by briefly adding a zero-length
# bookmark at that position and then deleting it again. Word (and this
# runtime) splits the underlying run when a bookmark boundary falls inside
# it, and the split survives bookmark deletion.
function Split-RunAtPosition($pos) {
    $zeroRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TmpSplitBookmark", $zeroRange) | Out-Null
    $d.Bookmarks.Item("TmpSplitBookmark").Delete()
}

# ===================================================================
# <sign> → + | - | λ    ==>   three separate paragraphs
# ===================================================================
$pSign = $d.Paragraphs.Item(34)
$pSign.Range.Text = "<sign> → +"

$pSign.Range.InsertParagraphAfter()
$pSignB = $d.Paragraphs.Item(35)
$pSignB.Range.Text = "<sign> →  - "

$pSignB.Range.InsertParagraphAfter()
$pSignC = $d.Paragraphs.Item(36)
$pSignC.Range.Text = "<sign> → λ"

# ===================================================================
# <digit> → 0 | 1 | ... | 9   ==>   ten separate paragraphs
# ===================================================================
$pDigit0 = $d.Paragraphs.Item(37)
$pDigit0.Range.Text = "<digit> → 0"

$pDigit0.Range.InsertParagraphAfter()
$pDigit1 = $d.Paragraphs.Item(38)
$pDigit1.Range.Text = "<digit> → 1"
Split-RunAtPosition ($pDigit1.Range.Start + 9)

$pDigit1.Range.InsertParagraphAfter()
$pDigit2 = $d.Paragraphs.Item(39)
$pDigit2.Range.Text = "<digit> → 2"
Split-RunAtPosition ($pDigit2.Range.Start + 10)

$pDigit2.Range.InsertParagraphAfter()
$pDigit3 = $d.Paragraphs.Item(40)
$pDigit3.Range.Text = "<digit> → 3"
Split-RunAtPosition ($pDigit3.Range.Start + 10)

$pDigit3.Range.InsertParagraphAfter()
$pDigit4 = $d.Paragraphs.Item(41)
$pDigit4.Range.Text = "<digit> → 4"
Split-RunAtPosition ($pDigit4.Range.Start + 10)

$pDigit4.Range.InsertParagraphAfter()
$pDigit5 = $d.Paragraphs.Item(42)
$pDigit5.Range.Text = "<digit> → 5"
Split-RunAtPosition ($pDigit5.Range.Start + 10)

$pDigit5.Range.InsertParagraphAfter()
$pDigit6 = $d.Paragraphs.Item(43)
$pDigit6.Range.Text = "<digit> → 6"
Split-RunAtPosition ($pDigit6.Range.Start + 10)

$pDigit6.Range.InsertParagraphAfter()
$pDigit7 = $d.Paragraphs.Item(44)
$pDigit7.Range.Text = "<digit> → 7"
Split-RunAtPosition ($pDigit7.Range.Start + 9)

$pDigit7.Range.InsertParagraphAfter()
$pDigit8 = $d.Paragraphs.Item(45)
$pDigit8.Range.Text = "<digit> → 8"
Split-RunAtPosition ($pDigit8.Range.Start + 9)

$pDigit8.Range.InsertParagraphAfter()
$pDigit9 = $d.Paragraphs.Item(46)
$pDigit9.Range.Text = "<digit> → 9"
Split-RunAtPosition ($pDigit9.Range.Start + 10)

# ===================================================================
# <id> → a | b | c | d   ==>   four separate paragraphs
# (first one keeps a bookmark named __DdeLink__572_463136898 around the
#  leading "<id> → " run)
# ===================================================================
$pIdA = $d.Paragraphs.Item(47)
$pIdA.Range.Text = "<id> → a"
$idBookmarkRange = $d.Range($pIdA.Range.Start, $pIdA.Range.Start + 7)
$d.Bookmarks.Add("__DdeLink__572_463136898", $idBookmarkRange) | Out-Null

$pIdA.Range.InsertParagraphAfter()
$pIdB = $d.Paragraphs.Item(48)
$pIdB.Range.Text = "<id> → b "

$pIdB.Range.InsertParagraphAfter()
$pIdC = $d.Paragraphs.Item(49)
$pIdC.Range.Text = "<id> → c "

$pIdC.Range.InsertParagraphAfter()
$pIdD = $d.Paragraphs.Item(50)
$pIdD.Range.Text = "<id> → d"
